$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Data": add the 2023 and 2022 rows at the top of the series, and
# correct the 2021 value (21 -> 21.6). Inserting two rows pushes every
# existing year row (2020..2005) down by two, which is exactly what the
# published table now shows (2005..2023, most recent year first).
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Rows("2:3").Insert()

# New rows need to keep "Fecha" formatted as text (like every other year in
# the column) instead of being auto-detected as numbers.
$wsData.Range("A2:A3").NumberFormat = "@"

$wsData.Range("A2").Value = "2023"
$wsData.Range("B2").Value = 19.5

$wsData.Range("A3").Value = "2022"
$wsData.Range("B3").Value = 18.5

# Correct the 2021 figure (now on row 4 after the insert).
$wsData.Range("B4").Value = 21.6

# ---------------------------------------------------------------------------
# Sheet "Metadata": record the July 2025 update, inserted right after the
# "observaciones" row (pushes "cita"/source rows down by one).
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# The first metadata row held a stray blank cell (empty string) next to a
# single-space cell; normalise it to a space so the unused empty shared
# string can drop out of the table entirely.
$wsMeta.Range("A1").Value = " "

$wsMeta.Rows("9:9").Insert()
$wsMeta.Range("A9").Value = "actualizacion"
$wsMeta.Range("B9").Value = "Julio 2025"
